$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.217.63"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.439.93"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "410.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.86%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.749"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.52%  "
$ws.Range("E10").Value = "  +16.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "43.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.90%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "3.969.88"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.43%  "
$ws.Range("E16").Value = "  +61.82%  "
$ws.Range("D17").Value = "3.456.67"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +15.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.09%  "
$ws.Range("D20").Value = "62.276.15"
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "408.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +29.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "90.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.73%  "
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "33.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "44.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.58%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.94"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.52%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.171"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0505"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.43%  "
$ws.Range("E37").Value = "  +2.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.132"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  +6.61%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("D49").Value = "2.124.02"
$ws.Range("E49").Value = "  +0.27%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.128"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.36%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.15%  "
